# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Fri Jul 26 22:29:31 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.038.93"
$ws.Range("E2").Value = "  +3.34%  "
$ws.Range("D3").Value = "3.276.48"
$ws.Range("E3").Value = "  +3.38%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'582.66"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D6").Value = "'182.85"
$ws.Range("E6").Value = "  +6.55%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.605"
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").Value = "3.277.22"
$ws.Range("E9").Value = "  +3.63%  "
$ws.Range("D10").Value = "'0.135"
$ws.Range("E10").Value = "  +8.13%  "
$ws.Range("D11").Value = "'6.74"
$ws.Range("E11").Value = "  +2.28%  "
$ws.Range("D12").Value = "'0.418"
$ws.Range("E12").Value = "  +6.54%  "
$ws.Range("D13").Value = "3.843.41"
$ws.Range("E13").Value = "  +3.48%  "
$ws.Range("D14").Value = "'0.138"
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").Value = "'28.58"
$ws.Range("E15").Value = "  +4.42%  "
$ws.Range("D16").Value = "68.014.34"
$ws.Range("E16").Value = "  +3.44%  "
$ws.Range("D17").Value = "'0.0000170"
$ws.Range("E17").Value = "  +3.74%  "
$ws.Range("D18").Value = "3.274.42"
$ws.Range("E18").Value = "  +3.25%  "
$ws.Range("D19").Value = "'5.85"
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("D20").Value = "'13.56"
$ws.Range("E20").Value = "  +5.20%  "
$ws.Range("D21").Value = "'377.69"
$ws.Range("E21").Value = "  +4.73%  "
$ws.Range("D22").Value = "'7.69"
$ws.Range("E22").Value = "  +5.70%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "'71.29"
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("D25").Value = "'0.515"
$ws.Range("E25").Value = "  +4.15%  "
$ws.Range("E26").Value = "  +5.87%  "
$ws.Range("D27").Value = "'9.68"
$ws.Range("E27").Value = "  -1.78%  "
$ws.Range("E28").Value = "  +2.33%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  +3.20%  "
$ws.Range("D31").Value = "'5.71"
$ws.Range("E31").Value = "  +6.28%  "
$ws.Range("D32").Value = "'22.93"
$ws.Range("E32").Value = "  +3.84%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.28"
$ws.Range("E33").Value = "  +6.76%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "'6.96"
$ws.Range("E35").Value = "  +5.13%  "
$ws.Range("D36").Value = "'1.53"
$ws.Range("E36").Value = "  +5.54%  "
$ws.Range("D37").Value = "'162.71"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("D38").Value = "'0.854"
$ws.Range("E38").Value = "  +2.26%  "
$ws.Range("E39").Value = "  +2.46%  "
$ws.Range("D40").Value = "'26.94"
$ws.Range("E40").Value = "  +1.95%  "
$ws.Range("E41").Value = "  +10.26%  "
$ws.Range("D42").Value = "'4.64"
$ws.Range("E42").Value = "  +10.78%  "
$ws.Range("D43").Value = "'2.64"
$ws.Range("E43").Value = "  +6.45%  "
$ws.Range("D44").Value = "'351.89"
$ws.Range("E44").Value = "  +7.07%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'25.79"
$ws.Range("E45").Value = "  +7.60%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.688.76"
$ws.Range("E46").Value = "  +1.79%  "
$ws.Range("D47").Value = "'40.86"
$ws.Range("D48").Value = "'0.0684"
$ws.Range("E48").Value = "  +3.83%  "
$ws.Range("D49").Value = "'0.0284"
$ws.Range("E49").Value = "  +3.71%  "
$ws.Range("E50").Value = "  +5.88%  "
$ws.Range("E51").Value = "  +1.70%  "
